$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BValues = @(1.127795056632408,1.003958681898894,0.9281657647420616,0.8973402314416035,0.8922253287704507,0.9277497956284719,1.085045507595737,1.39546760151859,1.62480816458384,1.729435608396159,1.769099300315929,1.760555069244958,1.732697888365408,1.715640247864712,1.617976580627726,1.558140093317888,1.523751864979772,1.512113418929118,1.564506873853531,1.740879031508371,1.856402529380205,1.794722038060172,1.561628413277447,1.311272081837501)
$CValues = @(0.3598284857417866,0.3154094122314461,0.2881830475973572,0.2770993894368701,0.2752596236199452,0.2880335244338426,0.3445026658286565,0.4556426078297022,0.5375934927687922,0.5749508005578718,0.5891089338602455,0.5860591951724246,0.5761153601371802,0.5700260197150442,0.5351537138422486,0.5137809386620233,0.5014951363717159,0.4973366024802885,0.5160553521383804,0.5790357823385648,0.6202657986888198,0.5982540807198689,0.5150270845917362,0.4255279095652327)
$DValues = @(0.01458388519450793,0.01486669361332105,0.01505050133698749,0.01512794859957367,0.01514096199419424,0.01505153553173066,0.01467928153634901,0.01403043402652315,0.01360389867316947,0.01342091067233042,0.01335322041944131,0.01336772717257695,0.01341530953409276,0.01344466435885039,0.01361608090221722,0.01372407888326421,0.01378723515194658,0.01380879672291968,0.01371247470682846,0.01340128980468869,0.01320726522790494,0.01330995892636189,0.0137177176331793,0.01419720625733234)
$FValues = @(0.9108577230097268,0.8883680273926302,0.8753927597837361,0.8703133503878888,0.8694824452255574,0.8753234163825283,0.9029293863305554,0.9637486917278011,1.01261141232483,1.035770637813968,1.044676113681334,1.04275210365094,1.03650056991907,1.03268903086709,1.011116784426648,0.9981226194824302,0.9907362807639544,0.9882503818812154,0.9994967963314139,1.038333104054004,1.06450591728327,1.05046406327682,0.9988752688927178,0.9465691796173274)
$GValues = @(0.00241298610574008,0.002417522126429627,0.002420453234484227,0.002421684513427946,0.002421891194199532,0.002420469690708626,0.002414519902827055,0.002404005021816189,0.002396974564431064,0.002393925425809822,0.002392792101908208,0.002393035237223819,0.002393831759959707,0.00239432242605525,0.002397176822176736,0.002398965995819072,0.002400009117538304,0.002400364714919787,0.002398774083525523,0.002393597224248301,0.002390338057625041,0.002392066208745988,0.002398860801955853,0.002406726992234598)
$IValues = @(0.5340523969013589,0.5337657693057949,0.5340845520731321,0.5343381574754744,0.5343877186806338,0.534087472384698,0.5338504438971157,0.537344668651663,0.54237277517678,0.5452049125037632,0.5463565423950953,0.5461049845342814,0.5452980665916058,0.5448141394537629,0.5421987212214887,0.5407343683415533,0.539943343825577,0.5396842909943729,0.540884942960929,0.5455329230153438,0.5490324401538942,0.5471221567805884,0.5408167098381611,0.5359706712129721)
$NValues = @(0.9464710712850817,0.9616761372101692,0.971518454279483,0.9756564190411972,0.9763511986539157,0.9715737457806917,0.9516085753968362,0.9164824562511491,0.8931397217436334,0.8830579758298249,0.8793177335469302,0.8801198104556711,0.8827487090835966,0.8843690860926081,0.8938094166353032,0.8997384832391475,0.9031992515717242,0.9043796777869844,0.8991020927521483,0.8819744323466487,0.8712323578136427,0.876924167376103,0.8993896427178427,0.9255530020422924)

$startRow = 2
$endRow = 25

for ($row = $startRow; $row -le $endRow; $row++) {
    $idx = $row - $startRow
    $ws.Cells.Item($row, "B").Value = $BValues[$idx]
    $ws.Cells.Item($row, "C").Value = $CValues[$idx]
    $ws.Cells.Item($row, "D").Value = $DValues[$idx]
    $ws.Cells.Item($row, "F").Value = $FValues[$idx]
    $ws.Cells.Item($row, "G").Value = $GValues[$idx]
    $ws.Cells.Item($row, "I").Value = $IValues[$idx]
    $ws.Cells.Item($row, "N").Value = $NValues[$idx]
}
